# Auto-generated edit script: updates Price (D), Volume(1h) (E), and Hora (G)
# columns for rows 2-51 on Sheet1, matching the "Updated symbol list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '286.45'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-9.87%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '17'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.98'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.40%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '17'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.034'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-4.11%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '17'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07275'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-5.99%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '17'
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.309'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.20%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '17'
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.522'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-11.05%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '17'
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9168'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '17'
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1197'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-5.35%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '17'
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1704'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-6.32%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '17'
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08639'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-6.20%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '17'
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04192'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.43%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '17'
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1052'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.08%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '17'
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001269'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.99%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '17'
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005965'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.39%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '17'
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.405'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.47%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '17'
# Row 17
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '17'
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.3282'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.09%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '17'
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.855'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.28%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '17'
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1343'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.57%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '17'
# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.42%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '17'
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.03850'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-4.09%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '17'
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001270'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.61%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '17'
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.003785'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-8.20%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '17'
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.15%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '17'
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003731'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '17'
# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '17'
# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '17'
# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '17'
# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '17'
# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '17'
# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '17'
# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '17'
# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '17'
# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '17'
# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '17'
# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '17'
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02308'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-9.34%'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '17'
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.04953'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-7.49%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '17'
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006365'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '221.85%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '17'
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007706'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.97%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '17'
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.96%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '17'
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007394'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '1.11%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '17'
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.006928'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-8.36%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '17'
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3091'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-9.92%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '17'
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006394'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.85%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '17'
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.36%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '17'
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '15.66%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '17'
# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.16%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '17'
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002105'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.36%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '17'
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002005'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.36%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '17'
